$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted as the most-recent entry
# for this market/product. It becomes the new row 520, which pushes the
# previously-existing rows 520-532 down to 521-533 (dimension grows from
# A1:R532 to A1:R533).
$ws.Rows.Item(520).Insert()

# Populate the newly inserted row 520 with the new record's data. The
# categorical columns (market, region, product, etc.) repeat the values
# already used throughout this sheet.
$ws.Cells.Item(520, 1).Value  = 10
$ws.Cells.Item(520, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(520, 3).Value  = "La Araucanía"
$ws.Cells.Item(520, 4).Value  = 45239
$ws.Cells.Item(520, 5).Value  = 9
$ws.Cells.Item(520, 6).Value  = 100112001
$ws.Cells.Item(520, 7).Value  = "Berenjena"
$ws.Cells.Item(520, 8).Value  = "Sin especificar"
$ws.Cells.Item(520, 9).Value  = "Primera"
$ws.Cells.Item(520, 10).Value = 200
$ws.Cells.Item(520, 11).Value = 13000
$ws.Cells.Item(520, 12).Value = 15000
$ws.Cells.Item(520, 13).Value = 14200
$ws.Cells.Item(520, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(520, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(520, 16).Value = 355
$ws.Cells.Item(520, 17).Value = 40
$ws.Cells.Item(520, 18).Value = "Hortaliza"
